# Auto-generated Excel COM-interop script to update market price data (columns H:N)
# across all 8 profession sheets, per the scheduled runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 1000  # H11: 753.75 -> 1000
$ws.Cells.Item(11, 9).Value = 1000  # I11: 753.75 -> 1000
$ws.Cells.Item(11, 11).Value = 1000  # K11: 753.75 -> 1000
$ws.Cells.Item(11, 13).Value = -860  # M11: -613.75 -> -860

$ws.Cells.Item(19, 8).Value = 846.44446  # H19: 797.8077 -> 846.44446
$ws.Cells.Item(19, 9).Value = 660  # I19: 519.1667 -> 660
$ws.Cells.Item(19, 10).Value = 1032.8889  # J19: 1036.6428 -> 1032.8889
$ws.Cells.Item(19, 11).Value = 660  # K19: 519.1667 -> 660
$ws.Cells.Item(19, 12).Value = 1032.8889  # L19: 1036.6428 -> 1032.8889
$ws.Cells.Item(19, 13).Value = -485  # M19: -344.1667 -> -485
$ws.Cells.Item(19, 14).Value = -1382.8889  # N19: -1386.6428 -> -1382.8889

$ws.Cells.Item(34, 8).Value = 4496  # H34: 17999.5 -> 4496
$ws.Cells.Item(34, 9).Value = 4496  # I34: 17999.5 -> 4496
$ws.Cells.Item(34, 11).Value = 4496  # K34: 17999.5 -> 4496
$ws.Cells.Item(34, 13).Value = -4293  # M34: -17796.5 -> -4293

$ws.Cells.Item(36, 8).Value = 4496  # H36: 17999.5 -> 4496
$ws.Cells.Item(36, 9).Value = 4496  # I36: 17999.5 -> 4496
$ws.Cells.Item(36, 11).Value = 4496  # K36: 17999.5 -> 4496
$ws.Cells.Item(36, 13).Value = -3781  # M36: -17284.5 -> -3781

$ws.Cells.Item(51, 8).Value = 2688.875  # H51: 2557.7083 -> 2688.875
$ws.Cells.Item(51, 9).Value = 2500  # I51: 2470.5881 -> 2500
$ws.Cells.Item(51, 10).Value = 3066.625  # J51: 2769.2856 -> 3066.625
$ws.Cells.Item(51, 11).Value = 2500  # K51: 2470.5881 -> 2500
$ws.Cells.Item(51, 12).Value = 3066.625  # L51: 2769.2856 -> 3066.625
$ws.Cells.Item(51, 13).Value = -2016  # M51: -1986.5881 -> -2016
$ws.Cells.Item(51, 14).Value = -4034.625  # N51: -3737.2856 -> -4034.625

$ws.Cells.Item(61, 8).Value = 115  # H61: 118.333336 -> 115
$ws.Cells.Item(61, 9).Value = 115  # I61: 118.333336 -> 115
$ws.Cells.Item(61, 11).Value = 345  # K61: 355.000008 -> 345
$ws.Cells.Item(61, 13).Value = -173  # M61: -183.000008 -> -173

$ws.Cells.Item(100, 8).Value = 1271.6818  # H100: 1273.6364 -> 1271.6818
$ws.Cells.Item(100, 9).Value = 1405.1052  # I100: 1463.7222 -> 1405.1052
$ws.Cells.Item(100, 10).Value = 426.66666  # J100: 418.25 -> 426.66666
$ws.Cells.Item(100, 11).Value = 1405.1052  # K100: 1463.7222 -> 1405.1052
$ws.Cells.Item(100, 12).Value = 426.66666  # L100: 418.25 -> 426.66666
$ws.Cells.Item(100, 13).Value = -864.1052  # M100: -922.7221999999999 -> -864.1052
$ws.Cells.Item(100, 14).Value = -1508.66666  # N100: -1500.25 -> -1508.66666

$ws.Cells.Item(111, 8).Value = 10793.667  # H111: 12195.5 -> 10793.667
$ws.Cells.Item(111, 9).Value = 1487.5  # I111: 6999.4287 -> 1487.5
$ws.Cells.Item(111, 10).Value = 20099.834  # J111: 17391.572 -> 20099.834
$ws.Cells.Item(111, 11).Value = 4462.5  # K111: 20998.2861 -> 4462.5
$ws.Cells.Item(111, 12).Value = 60299.50199999999  # L111: 52174.716 -> 60299.50199999999
$ws.Cells.Item(111, 13).Value = -1395.5  # M111: -17931.2861 -> -1395.5
$ws.Cells.Item(111, 14).Value = -66433.50199999999  # N111: -58308.716 -> -66433.50199999999

$ws.Cells.Item(112, 8).Value = 29641.422  # H112: 27618.44 -> 29641.422
$ws.Cells.Item(112, 10).Value = 36112.29  # J112: 33101.85 -> 36112.29
$ws.Cells.Item(112, 12).Value = 108336.87  # L112: 99305.54999999999 -> 108336.87
$ws.Cells.Item(112, 14).Value = -110552.87  # N112: -101521.55 -> -110552.87

$ws.Cells.Item(130, 8).Value = 132333.2  # H130: 132333.8 -> 132333.2
$ws.Cells.Item(130, 10).Value = 132333.2  # J130: 132333.8 -> 132333.2
$ws.Cells.Item(130, 12).Value = 132333.2  # L130: 132333.8 -> 132333.2
$ws.Cells.Item(130, 14).Value = -142373.2  # N130: -142373.8 -> -142373.2

$ws.Cells.Item(137, 8).Value = 1917  # H137: 1866 -> 1917
$ws.Cells.Item(137, 9).Value = 1854.15  # I137: 1776.409 -> 1854.15
$ws.Cells.Item(137, 11).Value = 5562.450000000001  # K137: 5329.227000000001 -> 5562.450000000001
$ws.Cells.Item(137, 13).Value = -3012.450000000001  # M137: -2779.227000000001 -> -3012.450000000001

$ws.Cells.Item(141, 8).Value = 1068.091  # H141: 1027.9166 -> 1068.091
$ws.Cells.Item(141, 9).Value = 1068.091  # I141: 1134.9 -> 1068.091
$ws.Cells.Item(141, 10).Value = 0  # J141: 493 -> 0
$ws.Cells.Item(141, 11).Value = 3204.273  # K141: 3404.7 -> 3204.273
$ws.Cells.Item(141, 12).Value = 0  # L141: 1479 -> 0
$ws.Cells.Item(141, 13).Value = 1975.727  # M141: 1775.3 -> 1975.727
$ws.Cells.Item(141, 14).ClearContents()  # N141: remove (was -11839)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1203.1578  # H2: 1328.55 -> 1203.1578
$ws.Cells.Item(2, 9).Value = 1166.8823  # I2: 1233.5625 -> 1166.8823
$ws.Cells.Item(2, 10).Value = 1511.5  # J2: 1708.5 -> 1511.5
$ws.Cells.Item(2, 11).Value = 1166.8823  # K2: 1233.5625 -> 1166.8823
$ws.Cells.Item(2, 12).Value = 1511.5  # L2: 1708.5 -> 1511.5
$ws.Cells.Item(2, 13).Value = -1053.8823  # M2: -1120.5625 -> -1053.8823
$ws.Cells.Item(2, 14).Value = -1737.5  # N2: -1934.5 -> -1737.5

$ws.Cells.Item(35, 8).Value = 2997.5  # H35: 2983.3333 -> 2997.5
$ws.Cells.Item(35, 9).Value = 2983  # I35: 2983.3333 -> 2983
$ws.Cells.Item(35, 10).Value = 3041  # J35: 0 -> 3041
$ws.Cells.Item(35, 11).Value = 2983  # K35: 2983.3333 -> 2983
$ws.Cells.Item(35, 12).Value = 3041  # L35: 0 -> 3041
$ws.Cells.Item(35, 13).Value = -2577  # M35: -2577.3333 -> -2577
$ws.Cells.Item(35, 14).Value = -3853  # N35: None -> -3853

$ws.Cells.Item(36, 8).Value = 2881.75  # H36: 5006.5 -> 2881.75
$ws.Cells.Item(36, 9).Value = 842.3333  # I36: 3675.3333 -> 842.3333
$ws.Cells.Item(36, 11).Value = 842.3333  # K36: 3675.3333 -> 842.3333
$ws.Cells.Item(36, 13).Value = -496.3333  # M36: -3329.3333 -> -496.3333

$ws.Cells.Item(45, 8).Value = 50949.5  # H45: 11096.4 -> 50949.5
$ws.Cells.Item(45, 9).Value = 99999  # I45: 12107.111 -> 99999
$ws.Cells.Item(45, 10).Value = 1900  # J45: 2000 -> 1900
$ws.Cells.Item(45, 11).Value = 99999  # K45: 12107.111 -> 99999
$ws.Cells.Item(45, 12).Value = 1900  # L45: 2000 -> 1900
$ws.Cells.Item(45, 13).Value = -99622  # M45: -11730.111 -> -99622
$ws.Cells.Item(45, 14).Value = -2654  # N45: -2754 -> -2654

$ws.Cells.Item(74, 8).Value = 8176.893  # H74: 7938.552 -> 8176.893
$ws.Cells.Item(74, 9).Value = 1752.4166  # I74: 1732.92 -> 1752.4166
$ws.Cells.Item(74, 11).Value = 1752.4166  # K74: 1732.92 -> 1752.4166
$ws.Cells.Item(74, 13).Value = -878.4166  # M74: -858.9200000000001 -> -878.4166

$ws.Cells.Item(77, 8).Value = 8176.893  # H77: 7938.552 -> 8176.893
$ws.Cells.Item(77, 9).Value = 1752.4166  # I77: 1732.92 -> 1752.4166
$ws.Cells.Item(77, 11).Value = 8762.083000000001  # K77: 8664.6 -> 8762.083000000001
$ws.Cells.Item(77, 13).Value = -4394.083000000001  # M77: -4296.6 -> -4394.083000000001

$ws.Cells.Item(116, 8).Value = 1203.1578  # H116: 1328.55 -> 1203.1578
$ws.Cells.Item(116, 9).Value = 1166.8823  # I116: 1233.5625 -> 1166.8823
$ws.Cells.Item(116, 10).Value = 1511.5  # J116: 1708.5 -> 1511.5
$ws.Cells.Item(116, 11).Value = 1166.8823  # K116: 1233.5625 -> 1166.8823
$ws.Cells.Item(116, 12).Value = 1511.5  # L116: 1708.5 -> 1511.5
$ws.Cells.Item(116, 13).Value = 1127.1177  # M116: 1060.4375 -> 1127.1177
$ws.Cells.Item(116, 14).Value = -6099.5  # N116: -6296.5 -> -6099.5

$ws.Cells.Item(132, 8).Value = 2730.7104  # H132: 2827.1667 -> 2730.7104
$ws.Cells.Item(132, 9).Value = 2480.484  # I132: 2582.9656 -> 2480.484
$ws.Cells.Item(132, 11).Value = 7441.451999999999  # K132: 7748.8968 -> 7441.451999999999
$ws.Cells.Item(132, 13).Value = -4911.451999999999  # M132: -5218.8968 -> -4911.451999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1203.1578  # H3: 1328.55 -> 1203.1578
$ws.Cells.Item(3, 9).Value = 1166.8823  # I3: 1233.5625 -> 1166.8823
$ws.Cells.Item(3, 10).Value = 1511.5  # J3: 1708.5 -> 1511.5
$ws.Cells.Item(3, 11).Value = 1166.8823  # K3: 1233.5625 -> 1166.8823
$ws.Cells.Item(3, 12).Value = 1511.5  # L3: 1708.5 -> 1511.5
$ws.Cells.Item(3, 13).Value = -1052.8823  # M3: -1119.5625 -> -1052.8823
$ws.Cells.Item(3, 14).Value = -1739.5  # N3: -1936.5 -> -1739.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 44053.875  # H31: 42309.76 -> 44053.875
$ws.Cells.Item(31, 9).Value = 54580.316  # I31: 51898.75 -> 54580.316
$ws.Cells.Item(31, 10).Value = 4053.4  # J31: 3953.8 -> 4053.4
$ws.Cells.Item(31, 11).Value = 54580.316  # K31: 51898.75 -> 54580.316
$ws.Cells.Item(31, 12).Value = 4053.4  # L31: 3953.8 -> 4053.4
$ws.Cells.Item(31, 13).Value = -54285.316  # M31: -51603.75 -> -54285.316
$ws.Cells.Item(31, 14).Value = -4643.4  # N31: -4543.8 -> -4643.4

$ws.Cells.Item(34, 8).Value = 44053.875  # H34: 42309.76 -> 44053.875
$ws.Cells.Item(34, 9).Value = 54580.316  # I34: 51898.75 -> 54580.316
$ws.Cells.Item(34, 10).Value = 4053.4  # J34: 3953.8 -> 4053.4
$ws.Cells.Item(34, 11).Value = 54580.316  # K34: 51898.75 -> 54580.316
$ws.Cells.Item(34, 12).Value = 4053.4  # L34: 3953.8 -> 4053.4
$ws.Cells.Item(34, 13).Value = -54378.316  # M34: -51696.75 -> -54378.316
$ws.Cells.Item(34, 14).Value = -4457.4  # N34: -4357.8 -> -4457.4

$ws.Cells.Item(60, 8).Value = 15598.75  # H60: 17278.6 -> 15598.75
$ws.Cells.Item(60, 10).Value = 16951.5  # J60: 19300.334 -> 16951.5
$ws.Cells.Item(60, 12).Value = 16951.5  # L60: 19300.334 -> 16951.5
$ws.Cells.Item(60, 14).Value = -17973.5  # N60: -20322.334 -> -17973.5

$ws.Cells.Item(99, 8).Value = 2978.5557  # H99: 2985 -> 2978.5557
$ws.Cells.Item(99, 9).Value = 2718.3333  # I99: 2727.75 -> 2718.3333
$ws.Cells.Item(99, 10).Value = 3499  # J99: 3499.5 -> 3499
$ws.Cells.Item(99, 11).Value = 2718.3333  # K99: 2727.75 -> 2718.3333
$ws.Cells.Item(99, 12).Value = 3499  # L99: 3499.5 -> 3499
$ws.Cells.Item(99, 13).Value = -1220.3333  # M99: -1229.75 -> -1220.3333
$ws.Cells.Item(99, 14).Value = -6495  # N99: -6495.5 -> -6495

$ws.Cells.Item(126, 8).Value = 2978.5557  # H126: 2985 -> 2978.5557
$ws.Cells.Item(126, 9).Value = 2718.3333  # I126: 2727.75 -> 2718.3333
$ws.Cells.Item(126, 10).Value = 3499  # J126: 3499.5 -> 3499
$ws.Cells.Item(126, 11).Value = 8154.999899999999  # K126: 8183.25 -> 8154.999899999999
$ws.Cells.Item(126, 12).Value = 10497  # L126: 10498.5 -> 10497
$ws.Cells.Item(126, 13).Value = -5684.999899999999  # M126: -5713.25 -> -5684.999899999999
$ws.Cells.Item(126, 14).Value = -15437  # N126: -15438.5 -> -15437

$ws.Cells.Item(132, 8).Value = 2707.7556  # H132: 2760.4792 -> 2707.7556
$ws.Cells.Item(132, 9).Value = 2498.487  # I132: 2540.025 -> 2498.487
$ws.Cells.Item(132, 10).Value = 4068  # J132: 3862.75 -> 4068
$ws.Cells.Item(132, 11).Value = 7495.461  # K132: 7620.075000000001 -> 7495.461
$ws.Cells.Item(132, 12).Value = 12204  # L132: 11588.25 -> 12204
$ws.Cells.Item(132, 13).Value = -4965.461  # M132: -5090.075000000001 -> -4965.461
$ws.Cells.Item(132, 14).Value = -17264  # N132: -16648.25 -> -17264

$ws.Cells.Item(133, 8).Value = 60226  # H133: 57113 -> 60226
$ws.Cells.Item(133, 10).Value = 60226  # J133: 57113 -> 60226
$ws.Cells.Item(133, 12).Value = 60226  # L133: 57113 -> 60226
$ws.Cells.Item(133, 14).Value = -65286  # N133: -62173 -> -65286

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(110, 8).Value = 34999.5  # H110: 35000 -> 34999.5
$ws.Cells.Item(110, 9).Value = 34999.5  # I110: 35000 -> 34999.5
$ws.Cells.Item(110, 11).Value = 104998.5  # K110: 105000 -> 104998.5
$ws.Cells.Item(110, 13).Value = -100908.5  # M110: -100910 -> -100908.5

$ws.Cells.Item(113, 8).Value = 495.45  # H113: 546.45 -> 495.45
$ws.Cells.Item(113, 10).Value = 440.69232  # J113: 519.1539 -> 440.69232
$ws.Cells.Item(113, 12).Value = 1322.07696  # L113: 1557.4617 -> 1322.07696
$ws.Cells.Item(113, 14).Value = -5662.07696  # N113: -5897.4617 -> -5662.07696

$ws.Cells.Item(137, 8).Value = 3422.5334  # H137: 3550.8572 -> 3422.5334
$ws.Cells.Item(137, 9).Value = 3048.1428  # I137: 3165.8572 -> 3048.1428
$ws.Cells.Item(137, 10).Value = 3750.125  # J137: 3935.8572 -> 3750.125
$ws.Cells.Item(137, 11).Value = 9144.428400000001  # K137: 9497.571599999999 -> 9144.428400000001
$ws.Cells.Item(137, 12).Value = 11250.375  # L137: 11807.5716 -> 11250.375
$ws.Cells.Item(137, 13).Value = -4044.428400000001  # M137: -4397.571599999999 -> -4044.428400000001
$ws.Cells.Item(137, 14).Value = -21450.375  # N137: -22007.5716 -> -21450.375

$ws.Cells.Item(140, 8).Value = 4216  # H140: 2843.3333 -> 4216
$ws.Cells.Item(140, 9).Value = 4288.1665  # I140: 2612 -> 4288.1665
$ws.Cells.Item(140, 10).Value = 3999.5  # J140: 4000 -> 3999.5
$ws.Cells.Item(140, 11).Value = 12864.4995  # K140: 7836 -> 12864.4995
$ws.Cells.Item(140, 12).Value = 11998.5  # L140: 12000 -> 11998.5
$ws.Cells.Item(140, 13).Value = -7684.499500000002  # M140: -2656 -> -7684.499500000002
$ws.Cells.Item(140, 14).Value = -22358.5  # N140: -22360 -> -22358.5

$ws.Cells.Item(141, 8).Value = 70295.81  # H141: 74782.2 -> 70295.81
$ws.Cells.Item(141, 9).Value = 11783.333  # I141: 13540 -> 11783.333
$ws.Cells.Item(141, 11).Value = 35349.999  # K141: 40620 -> 35349.999
$ws.Cells.Item(141, 13).Value = -30169.999  # M141: -35440 -> -30169.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3388.0715  # H132: 3649.92 -> 3388.0715
$ws.Cells.Item(132, 9).Value = 2617.5652  # I132: 2829.3 -> 2617.5652
$ws.Cells.Item(132, 11).Value = 7852.6956  # K132: 8487.900000000001 -> 7852.6956
$ws.Cells.Item(132, 13).Value = -5322.6956  # M132: -5957.900000000001 -> -5322.6956

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2013.2858  # H16: 2139.818 -> 2013.2858
$ws.Cells.Item(16, 9).Value = 2013.2858  # I16: 2139.818 -> 2013.2858
$ws.Cells.Item(16, 11).Value = 2013.2858  # K16: 2139.818 -> 2013.2858
$ws.Cells.Item(16, 13).Value = -1843.2858  # M16: -1969.818 -> -1843.2858

$ws.Cells.Item(55, 8).Value = 174.90909  # H55: 199.44444 -> 174.90909
$ws.Cells.Item(55, 9).Value = 167.4  # I55: 193.125 -> 167.4
$ws.Cells.Item(55, 11).Value = 167.4  # K55: 193.125 -> 167.4
$ws.Cells.Item(55, 13).Value = 5.599999999999994  # M55: -20.125 -> 5.599999999999994

$ws.Cells.Item(68, 8).Value = 2639  # H68: 2699.9 -> 2639
$ws.Cells.Item(68, 9).Value = 2639  # I68: 2699.9 -> 2639
$ws.Cells.Item(68, 11).Value = 2639  # K68: 2699.9 -> 2639
$ws.Cells.Item(68, 13).Value = -1890  # M68: -1950.9 -> -1890

$ws.Cells.Item(71, 8).Value = 2639  # H71: 2699.9 -> 2639
$ws.Cells.Item(71, 9).Value = 2639  # I71: 2699.9 -> 2639
$ws.Cells.Item(71, 11).Value = 13195  # K71: 13499.5 -> 13195
$ws.Cells.Item(71, 13).Value = -9451  # M71: -9755.5 -> -9451

$ws.Cells.Item(132, 8).Value = 4505.3  # H132: 4007.1333 -> 4505.3
$ws.Cells.Item(132, 9).Value = 4299.875  # I132: 3983.25 -> 4299.875
$ws.Cells.Item(132, 10).Value = 5327  # J132: 4102.6665 -> 5327
$ws.Cells.Item(132, 11).Value = 12899.625  # K132: 11949.75 -> 12899.625
$ws.Cells.Item(132, 12).Value = 15981  # L132: 12307.9995 -> 15981
$ws.Cells.Item(132, 13).Value = -10369.625  # M132: -9419.75 -> -10369.625
$ws.Cells.Item(132, 14).Value = -21041  # N132: -17367.9995 -> -21041

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 44237.168  # H69: 44046.145 -> 44237.168
$ws.Cells.Item(69, 9).Value = 43106  # I69: 43064.8 -> 43106
$ws.Cells.Item(69, 11).Value = 43106  # K69: 43064.8 -> 43106
$ws.Cells.Item(69, 13).Value = -42357  # M69: -42315.8 -> -42357

$ws.Cells.Item(72, 8).Value = 44237.168  # H72: 44046.145 -> 44237.168
$ws.Cells.Item(72, 9).Value = 43106  # I72: 43064.8 -> 43106
$ws.Cells.Item(72, 11).Value = 129318  # K72: 129194.4 -> 129318
$ws.Cells.Item(72, 13).Value = -125574  # M72: -125450.4 -> -125574
